# Bug fix for one trial with wrong date
# Updates median (B) and iqr (C) values for several institutions whose
# underlying computation changed once the trial with the incorrect date
# was corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Aalborg University Hospital
$ws.Range("B3").Value = 652.0000000000001
$ws.Range("C3").Value = 712.0000000000001

# Row 4: Aarhus University
$ws.Range("B4").Value = 565
$ws.Range("C4").Value = 821

# Row 5: Aarhus University Hospital
$ws.Range("C5").Value = 884.0000000000001

# Row 7: Bispebjerg and Frederiksberg Hospital
$ws.Range("B7").Value = 623.9999999999999

# Row 8: Copenhagen University Hospital
$ws.Range("C8").Value = 1052.5

# Row 15: Hvidovre and Amager Hospital
$ws.Range("C15").Value = 847.9999999999999

# Row 25: Norwegian University of Science and Technology
$ws.Range("C25").Value = 1456.5

# Row 36: Tampere University Hospital
$ws.Range("B36").Value = 381.9999999999999

# Row 44: University of Copenhagen
$ws.Range("B44").Value = 819

# Row 55: Zealand University Hospital
$ws.Range("B55").Value = 594
$ws.Range("C55").Value = 529.4999999999999

# Row 56: Denmark
$ws.Range("B56").Value = 617
$ws.Range("C56").Value = 885.0000000000001

# Row 57: Finland
$ws.Range("B57").Value = 702.5
$ws.Range("C57").Value = 1505

# Row 61: Total
$ws.Range("B61").Value = 690
$ws.Range("C61").Value = 1103
